{"js": "// Locate the paragraph that still has the original wording, insert a new\n// paragraph above it introducing the chocolatey pin command, and append a\n// colon to the original sentence.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst targetText = \"Try to use 6.0.x until the two applications are compatible. Here are the direct download links\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === targetText) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the 'Try to use 6.0.x...' paragraph\");\n}\n\n// Insert the new \"chocolatey pin\" paragraph right before the target\n// paragraph; it inherits the BodyText style from the paragraph it is\n// inserted next to.\nconst newPara = target.insertParagraph(\n  \"In chocolatey you can pin the virtualbox version with this syntax (it needs to be on a line by itself)\",\n  Word.InsertLocation.before\n);\nawait context.sync();\n\n// Append a separating space, then the verbatim choco command styled with\n// the existing \"Verbatim Char\" character style.\nnewPara.insertText(\" \", Word.InsertLocation.end);\nawait context.sync();\n\nconst codeRange = newPara.insertText(\n  \"choco install virtualbox --version 6.0.16\",\n  Word.InsertLocation.end\n);\ncodeRange.style = \"Verbatim Char\";\nawait context.sync();\n\n// The original sentence now gets a trailing colon.\ntarget.insertText(\":\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Insert a new \"chocolatey pin\" paragraph above the \"Try to use 6.0.x...\"\n# sentence, and append a trailing colon to that sentence.\n\n$d = $word.ActiveDocument\n\n$originalText = \"Try to use 6.0.x until the two applications are compatible. Here are the direct download links\"\n$newParaText = \"In chocolatey you can pin the virtualbox version with this syntax (it needs to be on a line by itself)\"\n$codeText = \"choco install virtualbox --version 6.0.16\"\n\nfunction Find-ParagraphByExactText($doc, $text) {\n  foreach ($p in $doc.Paragraphs) {\n    if ($p.Range.Text -eq ($text + \"`r\")) {\n      return $p\n    }\n  }\n  return $null\n}\n\n# 1) Find the paragraph that still has the original (un-colon-terminated) text\n#    and insert a brand-new empty paragraph right before it; the new\n#    paragraph inherits the BodyText style from its neighbor.\n$target = Find-ParagraphByExactText $d $originalText\n$target.Range.InsertParagraphBefore()\n\n# 2) Re-locate the (unchanged) target paragraph and grab the now-empty\n#    paragraph immediately preceding it -- that is our new paragraph.\n$target = Find-ParagraphByExactText $d $originalText\n$newPara = $target.Previous()\n$newPara.Range.Text = $newParaText\n\n# 3) Append \" choco install virtualbox --version 6.0.16\" right after the\n#    sentence we just wrote (still inside the same paragraph, before its\n#    paragraph mark).\n$newParaEnd = $newPara.Range.End - 1\n$insertionPoint = $d.Range($newParaEnd, $newParaEnd)\n$insertionPoint.InsertAfter(\" \" + $codeText)\n\n# 4) Style just the \"choco install ...\" run with the existing Verbatim\n#    character style used elsewhere in this document.\n$codeRange = $d.Content\n$codeRange.Find.ClearFormatting()\n$codeRange.Find.Text = $codeText\n$codeRange.Find.Execute() | Out-Null\n$codeRange.Style = \"Verbatim Char\"\n\n# 5) Add the trailing colon to the original sentence.\n$target = Find-ParagraphByExactText $d $originalText\n$targetEnd = $target.Range.End - 1\n$colonPoint = $d.Range($targetEnd, $targetEnd)\n$colonPoint.InsertAfter(\":\")\n"}
